$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new entry row (row 46) to the effort log, mirroring the
# existing rows: A=Date, B=Effort [h], D=description (shared string)
$ws.Range("A46").Value = 41234
$ws.Range("A46").NumberFormat = 'ddd\ dd/mm/yyyy'
$ws.Range("B46").Value = 1
$ws.Range("D46").Value = "Minor changes on documentation and setup"

# Keep the selection/active cell in sync with the newly added last row
$ws.Range("D46").Select()
